$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: 12062013_validation_10sect run, commented "test" ---
$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 41438
$ws.Range("C14").Value = "run_20130613T132541"
$ws.Range("B14").Value = "12062013_validation_10sect"
$ws.Range("E14").Value = "test"

# --- Row 15: 13062013_validation_10_sect_test run (filename shown in bold) ---
$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 41438
$ws.Range("C15").Value = "run_20130613T133513"
$ws.Range("B15").Value = "13062013_validation_10_sect_test"
$ws.Range("B15").Font.Bold = $true

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to A16 after the newly entered rows ---
$ws.Range("A16").Select()
